# The sheet "Artfynd" holds one observation record per row (rows 2..38,
# columns A..AY). The edit rotates every record down by one row: the old
# row 2 becomes row 3, old row 3 becomes row 4, ... old row 37 becomes
# row 38, and old row 38 wraps around to become the new row 2.
#
# We do this with full-row Range.Copy() operations (not .Value/.Value2
# assignment) so that text that looks like a date ("2018-10-31", "00:00",
# ...) is carried over as-is instead of being re-interpreted by Excel's
# smart cell-input parser and turned into a date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 38
$lastCol = "AY"

# 1) Stage a verbatim copy of the whole data block somewhere well outside
#    the used range, so we have a stable source to copy from even after
#    the original block below gets cleared.
$stageFirstRow = 100
$rowCount = $lastRow - $firstRow + 1
$stageLastRow = $stageFirstRow + $rowCount - 1

$ws.Range("A" + $firstRow + ":" + $lastCol + $lastRow).Copy($ws.Range("A" + $stageFirstRow))

# 2) Wipe the original block. Range.Copy() only touches destination cells
#    that have a source counterpart, so any cell that must end up blank
#    (because the row landing on top of it used to be shorter) needs an
#    explicit clear first.
$ws.Range("A" + $firstRow + ":" + $lastCol + $lastRow).ClearContents()

# 3) Copy staged rows back in rotated order: new row N (N = 3..lastRow)
#    gets old row N-1's data; new row firstRow gets old row lastRow's data.
for ($r = $firstRow + 1; $r -le $lastRow; $r++) {
    $oldRow = $r - 1
    $stageRow = $stageFirstRow + ($oldRow - $firstRow)
    $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow).Copy($ws.Range("A" + $r))
}

$stageRowForWrap = $stageFirstRow + ($lastRow - $firstRow)
$ws.Range("A" + $stageRowForWrap + ":" + $lastCol + $stageRowForWrap).Copy($ws.Range("A" + $firstRow))

# 4) Remove the staging area.
$ws.Range("A" + $stageFirstRow + ":" + $lastCol + $stageLastRow).ClearContents()
